# Update cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values, per the GitHub Actions refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.905.35'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.741.36'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5262'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2771'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06154'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("D11").Value = '1.736.09'
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07114'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6465'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.530'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9987'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9994'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '25.871.86'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006680'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").Value = '1.959.53'
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("E23").Value = '  +4.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.801'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.174'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.800'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.740'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.588'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.613'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9751'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6207'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.686'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.909'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9990'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3867'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7297'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.013'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05328'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1121'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.247'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.671'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.75%  '
